$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 149.83333
$ws.Range("I12").Value = 99.75
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 99.75
$ws.Range("L12").Value = 250
$ws.Range("M12").Value = 70.25
$ws.Range("N12").Value = -590

$ws.Range("H17").Value = 1940
$ws.Range("J17").Value = 1940
$ws.Range("L17").Value = 5820
$ws.Range("N17").Value = -6156

$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5936

$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5468

$ws.Range("H32").Value = 10099.7
$ws.Range("I32").Value = 8499.5
$ws.Range("J32").Value = 11166.5
$ws.Range("K32").Value = 8499.5
$ws.Range("L32").Value = 11166.5
$ws.Range("M32").Value = -8173.5
$ws.Range("N32").Value = -11818.5

$ws.Range("H33").Value = 372.89474
$ws.Range("I33").Value = 113.82353
$ws.Range("K33").Value = 113.82353
$ws.Range("M33").Value = 115.17647

$ws.Range("H40").Value = 6851.143
$ws.Range("J40").Value = 7224.231
$ws.Range("L40").Value = 7224.231
$ws.Range("N40").Value = -7574.231

$ws.Range("H41").Value = 4618.3335
$ws.Range("I41").Value = 3914.2
$ws.Range("K41").Value = 3914.2
$ws.Range("M41").Value = -3474.2

$ws.Range("H43").Value = 5750
$ws.Range("I43").Value = 5750
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 5750
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -5681
$ws.Range("N43").ClearContents()

$ws.Range("H51").Value = 23000.2
$ws.Range("I51").Value = 21250
$ws.Range("K51").Value = 21250
$ws.Range("M51").Value = -20766

$ws.Range("H69").Value = 5995
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 5995
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H96").Value = 3076.6667
$ws.Range("I96").Value = 2836.25
$ws.Range("K96").Value = 8508.75
$ws.Range("M96").Value = -7135.75

$ws.Range("H98").Value = 1768.7778
$ws.Range("I98").Value = 1433.625
$ws.Range("K98").Value = 1433.625
$ws.Range("M98").Value = 64.375

$ws.Range("H100").Value = 2118.1667
$ws.Range("I100").Value = 2118.1667
$ws.Range("K100").Value = 2118.1667
$ws.Range("M100").Value = -1577.1667

$ws.Range("H122").Value = 1768.7778
$ws.Range("I122").Value = 1433.625
$ws.Range("K122").Value = 4300.875
$ws.Range("M122").Value = -1850.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 37366.668
$ws.Range("I28").Value = 38550
$ws.Range("J28").Value = 35000
$ws.Range("K28").Value = 38550
$ws.Range("L28").Value = 35000
$ws.Range("M28").Value = -38358
$ws.Range("N28").Value = -35384

$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 5000
$ws.Range("N42").Value = -5972

$ws.Range("H92").Value = 550
$ws.Range("J92").Value = 550
$ws.Range("L92").Value = 550
$ws.Range("N92").Value = -5542

$ws.Range("H99").Value = 37366.668
$ws.Range("I99").Value = 38550
$ws.Range("J99").Value = 35000
$ws.Range("K99").Value = 38550
$ws.Range("L99").Value = 35000
$ws.Range("M99").Value = -35555
$ws.Range("N99").Value = -40990

$ws.Range("H110").Value = 837.6667
$ws.Range("I110").Value = 850
$ws.Range("K110").Value = 850
$ws.Range("M110").Value = 1195

$ws.Range("H133").Value = 99900
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99900
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99900
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -104960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1278.3572
$ws.Range("I20").Value = 1110.8889
$ws.Range("J20").Value = 1579.8
$ws.Range("K20").Value = 1110.8889
$ws.Range("L20").Value = 1579.8
$ws.Range("M20").Value = -863.8888999999999
$ws.Range("N20").Value = -2073.8

$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H80").Value = 294.9091
$ws.Range("I80").Value = 314
$ws.Range("J80").Value = 272
$ws.Range("K80").Value = 314
$ws.Range("L80").Value = 272
$ws.Range("M80").Value = 684
$ws.Range("N80").Value = -2268

$ws.Range("H83").Value = 294.9091
$ws.Range("I83").Value = 314
$ws.Range("J83").Value = 272
$ws.Range("K83").Value = 1570
$ws.Range("L83").Value = 1360
$ws.Range("M83").Value = 3422
$ws.Range("N83").Value = -11344

$ws.Range("H95").Value = 23541.334
$ws.Range("J95").Value = 23541.334
$ws.Range("L95").Value = 23541.334
$ws.Range("N95").Value = -29033.334

$ws.Range("H105").Value = 8000
$ws.Range("I105").Value = 8000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 8000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -6253
$ws.Range("N105").ClearContents()

$ws.Range("H140").Value = 100780
$ws.Range("J140").Value = 100780
$ws.Range("L140").Value = 100780
$ws.Range("N140").Value = -111140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 237.75
$ws.Range("I7").Value = 290.36365
$ws.Range("J7").Value = 122
$ws.Range("K7").Value = 290.36365
$ws.Range("L7").Value = 122
$ws.Range("M7").Value = -177.36365
$ws.Range("N7").Value = -348

$ws.Range("H16").Value = 1877.8
$ws.Range("I16").Value = 1133
$ws.Range("J16").Value = 2995
$ws.Range("K16").Value = 1133
$ws.Range("L16").Value = 2995
$ws.Range("M16").Value = -846
$ws.Range("N16").Value = -3569

$ws.Range("H22").Value = 755.1429000000001
$ws.Range("I22").Value = 753.82355
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 753.82355
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -403.82355
$ws.Range("N22").Value = -1500

$ws.Range("H105").Value = 5656.091
$ws.Range("I105").Value = 1801.75
$ws.Range("J105").Value = 7858.5713
$ws.Range("K105").Value = 1801.75
$ws.Range("L105").Value = 7858.5713
$ws.Range("M105").Value = -54.75
$ws.Range("N105").Value = -11352.5713

$ws.Range("H113").Value = 1877.8
$ws.Range("I113").Value = 1133
$ws.Range("J113").Value = 2995
$ws.Range("K113").Value = 1133
$ws.Range("L113").Value = 2995
$ws.Range("M113").Value = 1037
$ws.Range("N113").Value = -7335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1656.6666
$ws.Range("I8").Value = 1656.6666
$ws.Range("K8").Value = 4969.9998
$ws.Range("M8").Value = -4830.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3420
$ws.Range("I80").Value = 2699.6667
$ws.Range("J80").Value = 4500.5
$ws.Range("K80").Value = 2699.6667
$ws.Range("L80").Value = 4500.5
$ws.Range("M80").Value = -1701.6667
$ws.Range("N80").Value = -6496.5

$ws.Range("H83").Value = 3420
$ws.Range("I83").Value = 2699.6667
$ws.Range("J83").Value = 4500.5
$ws.Range("K83").Value = 13498.3335
$ws.Range("L83").Value = 22502.5
$ws.Range("M83").Value = -8506.333500000001
$ws.Range("N83").Value = -32486.5

$ws.Range("H123").Value = 49384
$ws.Range("J123").Value = 49384
$ws.Range("L123").Value = 49384
$ws.Range("N123").Value = -54284

$ws.Range("H132").Value = 6600.7144
$ws.Range("I132").Value = 5242.2
$ws.Range("K132").Value = 15726.6
$ws.Range("M132").Value = -13196.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 99995
$ws.Range("J36").Value = 99995
$ws.Range("L36").Value = 99995
$ws.Range("N36").Value = -101119

$ws.Range("H38").Value = 21005.8
$ws.Range("J38").Value = 21666.334
$ws.Range("L38").Value = 21666.334
$ws.Range("N38").Value = -22486.334

$ws.Range("H45").Value = 10040.5
$ws.Range("I45").Value = 10040.5
$ws.Range("K45").Value = 10040.5
$ws.Range("M45").Value = -9633.5

$ws.Range("H122").Value = 6666.6665
$ws.Range("I122").Value = 6666.6665
$ws.Range("K122").Value = 19999.9995
$ws.Range("M122").Value = -17549.9995

$ws.Range("H132").Value = 5279.846
$ws.Range("I132").Value = 5293.222
$ws.Range("J132").Value = 5249.75
$ws.Range("K132").Value = 15879.666
$ws.Range("L132").Value = 15749.25
$ws.Range("M132").Value = -13349.666
$ws.Range("N132").Value = -20809.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H45").Value = 33250
$ws.Range("I45").Value = 33250
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 33250
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -32759
$ws.Range("N45").ClearContents()

$ws.Range("H103").Value = 69663
$ws.Range("J103").Value = 69663
$ws.Range("L103").Value = 69663
$ws.Range("N103").Value = -72007

$ws.Range("H107").Value = 1668.2858
$ws.Range("I107").Value = 779.6667
$ws.Range("J107").Value = 7000
$ws.Range("K107").Value = 2339.0001
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = -419.0001000000002
$ws.Range("N107").Value = -24840

$ws.Range("H110").Value = 45644
$ws.Range("J110").Value = 45644
$ws.Range("L110").Value = 45644
$ws.Range("N110").Value = -53824

$ws.Range("H113").Value = 530.3333
$ws.Range("J113").Value = 651
$ws.Range("L113").Value = 1953
$ws.Range("N113").Value = -6293

$ws.Range("H132").Value = 1075
$ws.Range("I132").Value = 1250
$ws.Range("K132").Value = 3750
$ws.Range("M132").Value = -1220
